$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.383
$ws.Range("E2").Value = 0.447
$ws.Range("F2").Value = 0.197
$ws.Range("I2").Value = 0.0005727480334068931
$ws.Range("J2").Value = 0.0004979334718328064
$ws.Range("K2").Value = 1411.6
$ws.Range("L2").Value = 0.2409901835253948
$ws.Range("M2").Value = 319.5784
$ws.Range("N2").Value = 0.04439760492352148
$ws.Range("O2").Value = 0.2263944460187022
$ws.Range("P2").Value = 310.5684
$ws.Range("Q2").Value = 0.04314588571984274
$ws.Range("R2").Value = 0.2200116180221026
$ws.Range("S2").Value = 9.01
$ws.Range("T2").Value = 0.0281933947976459
$ws.Range("U2").Value = 7437.9
$ws.Range("V2").Value = 1.033314346841528
$ws.Range("W2").Value = 0.2455892731122089
$ws.Range("X2").Value = 0.1097854740990652
$ws.Range("Y2").Value = 0.1358037990131437
$ws.Range("Z2").Value = 1.333921285916345
$ws.Range("AB2").Value = 0.1047469695529461
$ws.Range("AC2").Value = -0.1047469695529461
$ws.Range("AD2").Value = 3352
$ws.Range("AE2").Value = 7.725641971595618
$ws.Range("AF2").Value = 3359.725641971595
$ws.Range("AG2").Value = -4078.174358028404
$ws.Range("AH2").Value = 0.3182213607142115
$ws.Range("AI2").Value = 0.2994694634040878
$ws.Range("AJ2").Value = -1.307138318671997
$ws.Range("AK2").Value = -1.078589447465863
$ws.Range("AN2").Value = 684.0816326530612
$ws.Range("AP2").Value = -832.2804812302865

# Row 3
$ws.Range("B3").Value = "Banco Patagonia S.A. (BASE:BPAT)"
$ws.Range("D3").Value = 0.375
$ws.Range("E3").Value = 0.482
$ws.Range("I3").Value = 0.007927390372591865
$ws.Range("J3").Value = 0.006836087282338959
$ws.Range("K3").Value = 199.2
$ws.Range("L3").Value = 0.4706994328922495
$ws.Range("M3").Value = 103.5504
$ws.Range("N3").Value = 0.1724690206528981
$ws.Range("O3").Value = 0.5198313253012048
$ws.Range("P3").Value = 103.5504
$ws.Range("Q3").Value = 0.1724690206528981
$ws.Range("R3").Value = 0.5198313253012048
$ws.Range("U3").Value = 199.9
$ws.Range("V3").Value = 0.3329447035309794
$ws.Range("W3").Value = 0.4363636363636363
$ws.Range("X3").Value = 0.09493726603005601
$ws.Range("Y3").Value = 0.3414263703335804
$ws.Range("Z3").Value = 0.5898047055672375
$ws.Range("AA3").Value = 0.004031956446791866
$ws.Range("AB3").Value = 0.09992160665577395
$ws.Range("AC3").Value = -0.09588965020898207
$ws.Range("AD3").Value = 72.40000000000001
$ws.Range("AE3").Value = 7.725641971595618
$ws.Range("AF3").Value = 80.12564197159563
$ws.Range("AG3").Value = -119.7743580284044
$ws.Range("AH3").Value = 0.1177408124394231
$ws.Range("AI3").Value = 0.1209944721153404
$ws.Range("AJ3").Value = -0.2492050934633298
$ws.Range("AK3").Value = -0.2590692515293432
$ws.Range("AN3").Value = 14.77551020408163
$ws.Range("AP3").Value = -24.44374653640905

# Row 4
$ws.Range("B4").Value = "Banco Santander Río S.A. (BASE:BRIO)"
$ws.Range("D4").Value = 0.383
$ws.Range("E4").Value = 0.39
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 243.6
$ws.Range("L4").Value = 0.2693200663349917
$ws.Range("M4").Value = 9.01
$ws.Range("N4").Value = 0.008204334365325076
$ws.Range("O4").Value = 0.03698686371100164
$ws.Range("P4").Value = -0
$ws.Range("Q4").Value = -0
$ws.Range("R4").Value = -0
$ws.Range("S4").Value = 9.01
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 1664.4
$ws.Range("V4").Value = 1.515570934256055
$ws.Range("W4").Value = 0.2455892731122089
$ws.Range("X4").Value = 0.1015340330672809
$ws.Range("Y4").Value = 0.144055240044928
$ws.Range("Z4").Value = 3.100118932146982
$ws.Range("AA4").Value = 0
$ws.Range("AB4").Value = 0.0996617683549787
$ws.Range("AC4").Value = -0.0996617683549787
$ws.Range("AD4").Value = 278.7
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 278.7
$ws.Range("AG4").Value = -1385.7
$ws.Range("AH4").Value = 0.2024112135957586
$ws.Range("AI4").Value = 0.1646579227224388
$ws.Range("AJ4").Value = 4.819826086956522
$ws.Range("AK4").Value = -49.13829787234035
$ws.Range("F4").ClearContents()
$ws.Range("AN4").ClearContents()
$ws.Range("AP4").ClearContents()

# Row 5
$ws.Range("D5").Value = 0.509
$ws.Range("E5").Value = 0.655
$ws.Range("F5").Value = 0.197
$ws.Range("K5").Value = 590.2
$ws.Range("L5").Value = 0.4395620764131973
$ws.Range("M5").Value = 168.1622
$ws.Range("N5").Value = 0.09729356630409627
$ws.Range("O5").Value = 0.2849240935276178
$ws.Range("P5").Value = 168.1622
$ws.Range("Q5").Value = 0.09729356630409627
$ws.Range("R5").Value = 0.2849240935276178
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 2330.2
$ws.Range("V5").Value = 1.348183290904883
$ws.Range("W5").Value = 0.4449302676215605
$ws.Range("X5").Value = 0.1035631023408472
$ws.Range("Y5").Value = 0.3413671652807133
$ws.Range("Z5").Value = 5.246971473231734
$ws.Range("AB5").Value = 0.1010221954184001
$ws.Range("AC5").Value = -0.1010221954184001
$ws.Range("AD5").Value = 502.6
$ws.Range("AF5").Value = 502.6
$ws.Range("AG5").Value = -1827.6
$ws.Range("AH5").Value = 0.2252801434334379
$ws.Range("AI5").Value = 0.222310686482661
$ws.Range("AJ5").Value = 18.42338709677423
$ws.Range("AK5").Value = 26.33429394812685

# Row 6
$ws.Range("D6").Value = 0.355
$ws.Range("E6").Value = -0.0554
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 33.7
$ws.Range("L6").Value = 0.04795105293113262
$ws.Range("M6").Value = 33.0858
$ws.Range("N6").Value = 0.02917361784675073
$ws.Range("O6").Value = 0.9817744807121661
$ws.Range("P6").Value = 33.0858
$ws.Range("Q6").Value = 0.02917361784675073
$ws.Range("R6").Value = 0.9817744807121661
$ws.Range("U6").Value = 571.9
$ws.Range("V6").Value = 0.5042765188255004
$ws.Range("W6").Value = 0.0345889356461049
$ws.Range("X6").Value = 0.1097854740990652
$ws.Range("Y6").Value = -0.07519653845296029
$ws.Range("Z6").Value = 0.6823300970873786
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.1047469695529461
$ws.Range("AC6").Value = -0.1047469695529461
$ws.Range("AD6").Value = 458.5
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 458.5
$ws.Range("AG6").Value = -113.4
$ws.Range("AH6").Value = 0.2878940097953033
$ws.Range("AI6").Value = 0.2466114457831325
$ws.Range("AJ6").Value = -0.111100225335554
$ws.Range("AK6").Value = -0.08809135399673733
$ws.Range("AN6").ClearContents()
$ws.Range("AP6").ClearContents()

# Row 7
$ws.Range("D7").Value = 0.5329999999999999
$ws.Range("E7").Value = 0.369
$ws.Range("F7").Value = 0.146
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 251.1
$ws.Range("L7").Value = 0.1367796056215274
$ws.Range("M7").Value = -0
$ws.Range("N7").Value = -0
$ws.Range("O7").Value = -0
$ws.Range("P7").Value = -0
$ws.Range("Q7").Value = -0
$ws.Range("R7").Value = -0
$ws.Range("U7").Value = 2341.9
$ws.Range("V7").Value = 1.097628421447319
$ws.Range("W7").Value = 0.1713408393039918
$ws.Range("X7").Value = 0.122387115042063
$ws.Range("Y7").Value = 0.04895372426192883
$ws.Range("Z7").Value = 1.244947782449478
$ws.Range("AA7").Value = 0
$ws.Range("AB7").Value = 0.11070547415127
$ws.Range("AC7").Value = -0.11070547415127
$ws.Range("AD7").Value = 1353
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 1353
$ws.Range("AG7").Value = -988.9000000000001
$ws.Range("AH7").Value = 0.388057133023576
$ws.Range("AI7").Value = 0.3937832882214267
$ws.Range("AJ7").Value = -0.8638944701668562
$ws.Range("AK7").Value = -0.903930530164534
$ws.Range("T7").ClearContents()
$ws.Range("AN7").ClearContents()
$ws.Range("AP7").ClearContents()

# Row 8
$ws.Range("B8").Value = "Grupo Supervielle S.A. (BASE:SUPV)"
$ws.Range("D8").Value = 0.535
$ws.Range("E8").Value = 0.526
$ws.Range("F8").Value = 1.168
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 43
$ws.Range("L8").Value = 0.08277189605389798
$ws.Range("M8").Value = 5.77
$ws.Range("N8").Value = 0.01699558173784978
$ws.Range("O8").Value = 0.1341860465116279
$ws.Range("P8").Value = 5.77
$ws.Range("Q8").Value = 0.01699558173784978
$ws.Range("R8").Value = 0.1341860465116279
$ws.Range("U8").Value = 182.3
$ws.Range("V8").Value = 0.5369661266568483
$ws.Range("W8").Value = 0.1231033495562554
$ws.Range("X8").Value = 0.1244072481628321
$ws.Range("Y8").Value = -0.001303898606576717
$ws.Range("Z8").Value = 9.293381037567091
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0.1115082097496771
$ws.Range("AC8").Value = -0.1115082097496771
$ws.Range("AD8").Value = 227.8
$ws.Range("AE8").Value = 0
$ws.Range("AF8").Value = 227.8
$ws.Range("AG8").Value = 45.5
$ws.Range("AH8").Value = 0.4015512074739997
$ws.Range("AI8").Value = 0.3530688158710477
$ws.Range("AJ8").Value = 0.1181818181818182
$ws.Range("AK8").Value = 0.09829336789803414
$ws.Range("AN8").ClearContents()
$ws.Range("AP8").ClearContents()

# Row 9
$ws.Range("D9").Value = 0.15
$ws.Range("E9").Value = 0.447
$ws.Range("K9").Value = 50.8
$ws.Range("L9").Value = 0.3937984496124031
$ws.Range("M9").Value = -0
$ws.Range("N9").Value = -0
$ws.Range("O9").Value = -0
$ws.Range("P9").Value = -0
$ws.Range("Q9").Value = -0
$ws.Range("R9").Value = -0
$ws.Range("U9").Value = 147.3
$ws.Range("V9").Value = 0.8987187309334961
$ws.Range("W9").Value = 0.2585241730279898
$ws.Range("X9").Value = 0.2411561166900341
$ws.Range("Y9").Value = 0.01736805633795566
$ws.Range("Z9").Value = 0.2281167108753316
$ws.Range("AB9").Value = 0.131456007662023
$ws.Range("AC9").Value = -0.131456007662023
$ws.Range("AD9").Value = 459
$ws.Range("AF9").Value = 459
$ws.Range("AG9").Value = 311.7
$ws.Range("AH9").Value = 0.7368759030341949
$ws.Range("AI9").Value = 0.6923076923076923
$ws.Range("AJ9").Value = 0.6553826745164003
$ws.Range("AK9").Value = 0.6044211751018033
$ws.Range("T9").ClearContents()
